$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fourth commit: add re-run feature for failed tc.
# Refreshed player low/high price data (re-run of the price fetch),
# which re-sorted the rows and updated the price columns. Values are
# stored as text (e.g. "29,900") to match the data source format, so we
# force each price cell to Text number format before writing the value
# and restore General afterwards (keeps the cell format identical to
# the original file while avoiding Excel's automatic number coercion).
$data = @(
    @{Row=2; A="Reis"; B="29,900"; C="33,000"},
    @{Row=3; A="Jony"; B="31,100"; C="34,400"},
    @{Row=4; A="Doğukan Sinik"; B="38,000"; C="42,100"},
    @{Row=5; A="Johan Caballero"; B="39,500"; C="43,600"},
    @{Row=6; A="Tobias Mohr"; B="38,800"; C="42,900"},
    @{Row=7; A="Emil Hansson"; B="30,400"; C="33,600"},
    @{Row=8; A="Song Min Kyu"; B="39,500"; C="43,600"},
    @{Row=9; A="Kylian Hazard"; B="29,000"; C="32,100"},
    @{Row=10; A="Octavian Popescu"; B="39,500"; C="43,600"},
    @{Row=11; A="Gustav Mendonca Wikheim"; B="32,700"; C="36,200"},
    @{Row=12; A="Renaldo Cephas"; B="39,500"; C="43,600"},
    @{Row=13; A="Mounir Chouiar"; B="0"; C="0"},
    @{Row=14; A="Michael Johnston"; B="0"; C="0"},
    @{Row=15; A="Michael Johnston"; B="38,300"; C="42,400"},
    @{Row=16; A="Rodrigo Martins"; B="39,500"; C="43,600"},
    @{Row=17; A="Camacho"; B="31,300"; C="34,600"},
    @{Row=18; A="Joaquín Valiente"; B="39,500"; C="43,600"},
    @{Row=19; A="Lameck Banda"; B="39,500"; C="43,600"},
    @{Row=20; A="Washington Corozo"; B="39,500"; C="43,600"},
    @{Row=21; A="Stipe Biuk"; B="39,500"; C="43,600"},
    @{Row=22; A="Léo Jabá"; B="39,500"; C="43,600"},
    @{Row=23; A="Carlos Forbs"; B="39,500"; C="43,600"}
)

foreach ($item in $data) {
    $row = $item.Row

    $ws.Cells.Item($row, 1).Value = $item.A

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 2).NumberFormat = "General"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 3).NumberFormat = "General"
}
